# Nhom 12 cap nhat
# Fill in the evaluation data for group 12 (rows 9-13) on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Student roster for group 12: MSSV (col C) + full name (col D)
$students = @(
    @{ Row = 9;  MSSV = "0712110"; Name = "Đoàn Nguyên Dương" },
    @{ Row = 10; MSSV = "0712127"; Name = "Trần Đức Duy" },
    @{ Row = 11; MSSV = "0712129"; Name = "Trịnh Khắc Duy" },
    @{ Row = 12; MSSV = "0712138"; Name = "Nguyễn Việt Hằng" },
    @{ Row = 13; MSSV = "0712325"; Name = "Trần Nam Phương" }
)

# Per-assignment completion percentages (columns E..K = Phan cong 1..7).
# Default is 100% (1); row 11 has two 75% (0.75) marks.
$scores = @{
    9  = @(1, 1, 1, 1, 1, 1, 1)
    10 = @(1, 1, 1, 1, 1, 1, 1)
    11 = @(1, 0.75, 1, 1, 1, 0.75, 1)
    12 = @(1, 1, 1, 1, 1, 1, 1)
    13 = @(1, 1, 1, 1, 1, 1, 1)
}

$cols = @("E", "F", "G", "H", "I", "J", "K")

# Write every MSSV first, then every name, so the shared-string table fills
# up in the same order (all IDs, then all names) as the source workbook.
foreach ($s in $students) {
    $ws.Range("C" + $s.Row).Value = $s.MSSV
}
foreach ($s in $students) {
    $ws.Range("D" + $s.Row).Value = $s.Name
}

foreach ($s in $students) {
    $r = $s.Row
    $vals = $scores[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $cell = $ws.Range($cols[$i] + $r)
        $cell.Value = $vals[$i]
        $cell.NumberFormat = "0%"
    }
}
